$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 477 - this shifts all existing rows
# (477..564) down to (479..566), and creates two fresh blank rows at
# 477 and 478 that we now need to populate.
$ws.Rows("477:478").Insert()

# New row 477 - "Primera" quality entry for Terminal Hortofrutícola Agro
# Chillán / Betarraga, dated 44995, origin changed to Provincia de Diguillín.
$ws.Cells.Item(477, 1).Value = 7
$ws.Cells.Item(477, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(477, 3).Value = "Ñuble"
$ws.Cells.Item(477, 4).Value = 44995
$ws.Cells.Item(477, 5).Value = 16
$ws.Cells.Item(477, 6).Value = 100114014
$ws.Cells.Item(477, 7).Value = "Betarraga"
$ws.Cells.Item(477, 8).Value = "Sin especificar"
$ws.Cells.Item(477, 9).Value = "Primera"
$ws.Cells.Item(477, 10).Value = 300
$ws.Cells.Item(477, 11).Value = 800
$ws.Cells.Item(477, 12).Value = 800
$ws.Cells.Item(477, 13).Value = 800
$ws.Cells.Item(477, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(477, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(477, 16).Value = 160
$ws.Cells.Item(477, 17).Value = 5
$ws.Cells.Item(477, 18).Value = "Hortaliza"

# New row 478 - "Segunda" quality entry, same market/date as row 477.
$ws.Cells.Item(478, 1).Value = 7
$ws.Cells.Item(478, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(478, 3).Value = "Ñuble"
$ws.Cells.Item(478, 4).Value = 44995
$ws.Cells.Item(478, 5).Value = 16
$ws.Cells.Item(478, 6).Value = 100114014
$ws.Cells.Item(478, 7).Value = "Betarraga"
$ws.Cells.Item(478, 8).Value = "Sin especificar"
$ws.Cells.Item(478, 9).Value = "Segunda"
$ws.Cells.Item(478, 10).Value = 200
$ws.Cells.Item(478, 11).Value = 600
$ws.Cells.Item(478, 12).Value = 600
$ws.Cells.Item(478, 13).Value = 600
$ws.Cells.Item(478, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(478, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(478, 16).Value = 120
$ws.Cells.Item(478, 17).Value = 5
$ws.Cells.Item(478, 18).Value = "Hortaliza"

# Preserve the date-style formatting (numFmt 165) on the new date cells -
# match the style already used by every other "Fecha" cell in column D.
$ws.Cells.Item(477, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(478, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
